# Plantilla4.xlsx edit: add a "Cantidad Matriculas" column, drop the old
# "Porcentaje" value, change the career from "Negocios" to "Sistemas" and
# add a second data row for year 2023.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the brand-new strings in the same order they first appear so the
# regenerated shared-string table lines up with the target workbook
# ("Sistemas" is introduced before "Cantidad Matriculas").
$ws.Range("B2").Value = "Sistemas"
$ws.Range("C1").Value = "Cantidad Matriculas"

# Re-assert the header row (columns shift right to make room for C1)
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Carrera"
$ws.Range("D1").Value = "Porcentaje"
$ws.Range("E1").Value = "Periodo"
$ws.Range("F1").Value = "Año"

# Row 2: Carrera -> Sistemas (already set above), drop the old Porcentaje
# value, keep Periodo/Año but shifted one column to the right (E/F).
$ws.Range("A2").Value = 1
$ws.Range("C2").Value = 12
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = "ENE-ABR"
$ws.Range("F2").Value = 2022

# New row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Sistemas"
$ws.Range("C3").Value = 12
$ws.Range("E3").Value = "ENE-ABR"
$ws.Range("F3").Value = 2023

# Widen the new "Cantidad Matriculas" column
$ws.Columns.Item(3).ColumnWidth = 19.25

# Move the active selection like the author left it
$ws.Range("E6").Select()
